# Spelontwerp - v0.1: fill in the "Architectuur van de game" section with the
# real write-up (replacing the placeholder "(Oa: ...)" bullet and the three
# trailing blank paragraphs with the finished text).

$d = $word.ActiveDocument

$apos = [char]0x2019

$content1 = "Voor elk apart Object, wordt een eigen Klasse geschreven in de programmeertaal Java. Zo zijn er klassen voor Unit, Field, Tower, enzovoorts. De klassen die met het veld te maken hebben, zoals Field, Tower en Tree (obstakel), hebben allemaal de klasse Field als basis (spreekt voor zich voor de klasse Field). Met andere woorden, klassen als Tower en Tree breiden de klasse Field uit. Zo hebben ze deze functies en eigenschappen, met nog wat extra functies en eigenschappen die van belang zijn voor dat type."

$content2 = "Zoals hierboven al even is genoemd, wordt het programma geschreven in de programmeertaal Java. Het programma moet in internet browsers werken, zodat het spel gespeeld kan worden op internet. Het spel is daarnaast niet geschikt voor multiplayer, maar er kan slechts tegen het spel worden gespeeld."

$content3 = "Er zijn een heleboel verschillende programma${apos}s waarin geprogrammeerd kan worden, als het Java betreft. Iedereen heeft zijn eigen voorkeur en daarom zijn er binnen onze groep in ieder geval 2 programma${apos}s in gebruik (Netbeans en Eclipse). Om deze programma${apos}s niet te laten communiceren met elkaar, wat betreft de zogenaamde projecten die beide programma${apos}s gebruiken, zetten we alleen de .java en de .class bestanden op de server. Zo kan er nooit een fout optreden in de programma${apos}s die we gebruiken, behalve natuurlijk programmeerfouten."

# 1) The placeholder paragraph "(Oa: software architectuur, ontwikkeltaal,
#    platform, frameworks, multiplayer etc.)" becomes empty, and the first
#    real paragraph of content is appended right after it (inheriting the
#    "Geenafstand"/No Spacing style of the placeholder paragraph).
$oldPlaceholder = "(Oa: software architectuur, ontwikkeltaal, platform, frameworks, multiplayer etc.)"
$d.Content.Find.Execute($oldPlaceholder, $false, $false, $false, $false, $false, $true, 1, $false, ("^p" + $content1), 2) | Out-Null

# 2) Locate the three originally-empty trailing paragraphs (right after the
#    paragraph that now holds $content1) and expand them into the
#    empty / content2 / empty / content3 / empty sequence the final doc has.
$anchor = $d.Content.Find.Execute($content1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterContent1 = $d.Paragraphs.Item($d.Range(0, $d.Content.End).Paragraphs.Count)

# Find the paragraph index of the paragraph containing $content1, then work
# relative to it so this is robust regardless of exact counts.
$idx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq $content1) {
        $idx = $i
        break
    }
}

$p1 = $d.Paragraphs.Item($idx + 1)
$p1.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($idx + 3)
$p3.Range.InsertParagraphAfter()

$d.Paragraphs.Item($idx + 2).Range.Text = $content2
$d.Paragraphs.Item($idx + 4).Range.Text = $content3

# 3) All five of these paragraphs (the blank separators and the two new
#    content paragraphs) use the "Geenafstand" (No Spacing) style, matching
#    the rest of this section.
for ($i = $idx + 1; $i -le $idx + 5; $i++) {
    $d.Paragraphs.Item($i).Style = "Geenafstand"
}
